$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text interpretation so that
# numeric-looking strings (e.g. "34.80", "0.0680") keep their exact
# formatting/trailing zeros instead of being coerced into real numbers.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = "42.450.43"
$ws.Range("E2").Value = "  -1.35%  "
$ws.Range("D3").Value = "2.283.90"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue "D5" "303.68"
$ws.Range("E5").Value = "  +0.98%  "
Set-TextValue "D6" "95.35"
$ws.Range("E6").Value = "  -3.03%  "
$ws.Range("E7").Value = "  -2.92%  "
$ws.Range("E8").Value = "  +0.04%  "
Set-TextValue "D9" "0.494"
$ws.Range("E9").Value = "  -3.24%  "
Set-TextValue "D10" "34.80"
$ws.Range("E10").Value = "  -4.16%  "
Set-TextValue "D11" "0.0779"
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("E12").Value = "  +1.60%  "
Set-TextValue "D13" "18.02"
$ws.Range("E13").Value = "  +0.58%  "
Set-TextValue "D14" "6.75"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").Value = "2.639.38"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").Value = "2.283.23"
$ws.Range("E16").Value = "  -0.69%  "
Set-TextValue "D17" "0.768"
$ws.Range("E17").Value = "  -1.81%  "
$ws.Range("D18").Value = "42.372.67"
Set-TextValue "D19" "12.75"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("E21").Value = "  -2.81%  "
Set-TextValue "D22" "67.03"
$ws.Range("E22").Value = "  -1.89%  "
Set-TextValue "D23" "235.78"
$ws.Range("E23").Value = "  -2.76%  "
Set-TextValue "D24" "2.13"
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -2.13%  "
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("E28").Value = "  +17.15%  "
Set-TextValue "D29" "166.71"
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").Value = "  -1.63%  "
Set-TextValue "D31" "32.34"
$ws.Range("E31").Value = "  -2.85%  "
$ws.Range("E32").Value = "  +0.05%  "
Set-TextValue "D33" "17.69"
$ws.Range("E33").Value = "  -0.38%  "
Set-TextValue "D34" "4.93"
$ws.Range("E34").Value = "  -1.91%  "
Set-TextValue "D35" "4.42"
$ws.Range("E35").Value = "  -7.58%  "
$ws.Range("E36").Value = "  -2.45%  "
Set-TextValue "D37" "0.0680"
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("E39").Value = "  -2.59%  "
$ws.Range("E40").Value = "  -2.33%  "
$ws.Range("E41").Value = "  -4.40%  "
$ws.Range("D42").Value = "1.985.58"
$ws.Range("E42").Value = "  -0.66%  "
$ws.Range("E43").Value = "  -3.78%  "
Set-TextValue "D44" "10.08"
$ws.Range("E44").Value = "  -1.41%  "
Set-TextValue "D45" "18.13"
$ws.Range("E45").Value = "  +3.75%  "
$ws.Range("E46").Value = "  -8.49%  "
$ws.Range("E47").Value = "  -2.29%  "
Set-TextValue "D48" "2.91"
$ws.Range("E48").Value = "  +4.52%  "
Set-TextValue "D49" "53.37"
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("D50").Value = "2.506.35"

# Row 51: coin changed from BitcoinSV to TrustWalletToken
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D51" "1.12"
$ws.Range("E51").Value = "  +0.20%  "
